$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows right before row 293. This shifts the existing rows
# 293:332 down to 298:337 (dimension grows from A1:R332 to A1:R337),
# matching the diff's row-shift pattern.
$ws.Range("293:297").EntireRow.Insert()

# Populate the 5 newly-inserted rows (293-297) with a new week of data,
# following the same constant-column layout as the rest of the table
# (A=Mercado ID, B=Mercado, C=Region, E=Codreg, F=Categoria ID,
#  G=Categoria, N=Unidad de comercializacion, Q=Kg o Unidades, R=Clasificacion).

# Row 293
$ws.Range("A293").Value = 9
$ws.Range("B293").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C293").Value = "Metropolitana"
$ws.Range("D293").Value = 44491
$ws.Range("E293").Value = 13
$ws.Range("F293").Value = 100112006
$ws.Range("G293").Value = "Repollo"
$ws.Range("H293").Value = "Crespo record"
$ws.Range("I293").Value = "Primera"
$ws.Range("J293").Value = 5200
$ws.Range("K293").Value = 650
$ws.Range("L293").Value = 700
$ws.Range("M293").Value = 675
$ws.Range("N293").Value = "$/unidad"
$ws.Range("O293").Value = "Región Metropolitana"
$ws.Range("P293").Value = 675
$ws.Range("Q293").Value = 1
$ws.Range("R293").Value = "Hortaliza"

# Row 294
$ws.Range("A294").Value = 9
$ws.Range("B294").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C294").Value = "Metropolitana"
$ws.Range("D294").Value = 44491
$ws.Range("E294").Value = 13
$ws.Range("F294").Value = 100112006
$ws.Range("G294").Value = "Repollo"
$ws.Range("H294").Value = "Crespo record"
$ws.Range("I294").Value = "Primera"
$ws.Range("J294").Value = 3400
$ws.Range("K294").Value = 650
$ws.Range("L294").Value = 700
$ws.Range("M294").Value = 675
$ws.Range("N294").Value = "$/unidad"
$ws.Range("O294").Value = "Región de O'Higgins"
$ws.Range("P294").Value = 675
$ws.Range("Q294").Value = 1
$ws.Range("R294").Value = "Hortaliza"

# Row 295
$ws.Range("A295").Value = 9
$ws.Range("B295").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C295").Value = "Metropolitana"
$ws.Range("D295").Value = 44491
$ws.Range("E295").Value = 13
$ws.Range("F295").Value = 100112006
$ws.Range("G295").Value = "Repollo"
$ws.Range("H295").Value = "Crespo record"
$ws.Range("I295").Value = "Segunda"
$ws.Range("J295").Value = 2500
$ws.Range("K295").Value = 500
$ws.Range("L295").Value = 550
$ws.Range("M295").Value = 525
$ws.Range("N295").Value = "$/unidad"
$ws.Range("O295").Value = "Región Metropolitana"
$ws.Range("P295").Value = 525
$ws.Range("Q295").Value = 1
$ws.Range("R295").Value = "Hortaliza"

# Row 296
$ws.Range("A296").Value = 9
$ws.Range("B296").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C296").Value = "Metropolitana"
$ws.Range("D296").Value = 44491
$ws.Range("E296").Value = 13
$ws.Range("F296").Value = 100112006
$ws.Range("G296").Value = "Repollo"
$ws.Range("H296").Value = "Crespo record"
$ws.Range("I296").Value = "Segunda"
$ws.Range("J296").Value = 1600
$ws.Range("K296").Value = 500
$ws.Range("L296").Value = 550
$ws.Range("M296").Value = 525
$ws.Range("N296").Value = "$/unidad"
$ws.Range("O296").Value = "Región de O'Higgins"
$ws.Range("P296").Value = 525
$ws.Range("Q296").Value = 1
$ws.Range("R296").Value = "Hortaliza"

# Row 297
$ws.Range("A297").Value = 9
$ws.Range("B297").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C297").Value = "Metropolitana"
$ws.Range("D297").Value = 44491
$ws.Range("E297").Value = 13
$ws.Range("F297").Value = 100112006
$ws.Range("G297").Value = "Repollo"
$ws.Range("H297").Value = "Morada(o)"
$ws.Range("I297").Value = "Primera"
$ws.Range("J297").Value = 2500
$ws.Range("K297").Value = 800
$ws.Range("L297").Value = 900
$ws.Range("M297").Value = 850
$ws.Range("N297").Value = "$/unidad"
$ws.Range("O297").Value = "Región de O'Higgins"
$ws.Range("P297").Value = 850
$ws.Range("Q297").Value = 1
$ws.Range("R297").Value = "Hortaliza"
